$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the actual time values for the "Seeders and migrations" row (row 17)
$ws.Range("B17").Value = 1.2
$ws.Range("C17").Value = 1.3

# Update the selected cell to match the edited row
$ws.Range("B17").Select()

$wb.Save()
